# "update to the new UI"
#
# The deck started as a single blank "Title and Content" slide. The
# target deck has 5 "Title and Content" slides, each with a title and a
# body paragraph (sized down to 18pt to match the deck's default text
# style), plus a trailing blank paragraph in the body placeholder.

$p = $ppt.ActivePresentation

# Drop the original placeholder slide - every slide in the final deck is
# freshly authored content, not an edit of the old single slide.
$p.Slides.Item(1).Delete()

$titles = @("Introduction", "Slide 2", "Slide 3", "Slide 4", "Slide 5")
$bodies = @("Introduction", "Parties Involved", "Work Description", "Agreement Terms", "Conclusion")

for ($i = 0; $i -lt 5; $i++) {
    # Layout index 2 == "Title and Content" (same layout the original
    # slide used: Title placeholder + idx=1 Content Placeholder).
    $slide = $p.Slides.Add($i + 1, 2)

    $titleShape = $slide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Text = $titles[$i]

    $bodyShape = $slide.Shapes.Item(2)
    $bodyTr = $bodyShape.TextFrame.TextRange
    # Trailing CR leaves a blank second paragraph in the body, matching
    # the target's empty <a:p/> after the content line.
    $bodyTr.Text = $bodies[$i] + "`r"
    $bodyTr.Characters(1, $bodies[$i].Length).Font.Size = 18
}
